$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F22").Value = 168
$ws.Range("G22").Value = 19365.36
$ws.Range("B41").Value = 82303.97
$ws.Range("F115").Value = 22
$ws.Range("G115").Value = 20172.68
$ws.Range("B116").Value = 29881.11
$ws.Range("F125").Value = 197
$ws.Range("G125").Value = 22099.46
$ws.Range("F136").Value = 127
$ws.Range("G136").Value = 5935.98
$ws.Range("B147").Value = 98754.44
$ws.Range("B151").Value = 64196
$ws.Range("F151").Value = 1
$ws.Range("G151").Value = 32143.58
$ws.Range("B152").Value = 65258
$ws.Range("F152").Value = 2
$ws.Range("G152").Value = 64287.16
$ws.Range("F182").Value = 4
$ws.Range("G182").Value = 262.68
$ws.Range("B184").Value = 29824.51
$ws.Range("F212").Value = 1
$ws.Range("G212").Value = 366.29
$ws.Range("B220").Value = 31762.04
$ws.Range("F237").Value = 92
$ws.Range("G237").Value = 7544
$ws.Range("F238").Value = 26
$ws.Range("G238").Value = 1602.9
$ws.Range("B241").Value = 64329
$ws.Range("E241").Value = 128.32
$ws.Range("F241").Value = 1
$ws.Range("G241").Value = 120.69
$ws.Range("B242").Value = 57552
$ws.Range("E242").Value = 136.86
$ws.Range("F242").Value = -5
$ws.Range("G242").Value = -603.45
$ws.Range("B250").Value = 101508.9
$ws.Range("F282").Value = 115
$ws.Range("G282").Value = 8544.5
$ws.Range("B283").Value = 115110.19
$ws.Range("F287").Value = 2265
$ws.Range("G287").Value = 41902.5
$ws.Range("B294").Value = 64125.6
$ws.Range("F311").Value = 11
$ws.Range("G311").Value = 7441.72
$ws.Range("B318").Value = 24906.55
$ws.Range("B370").Value = 66194
$ws.Range("C370").Value = 'HIM-Total Care Baby Pants Diapers-M-9s'
$ws.Range("F370").Value = 39
$ws.Range("G370").Value = 3341.52
$ws.Range("B371").Value = 64983
$ws.Range("C371").Value = 'HIM-TOTAL CARE BABY PANTS DIAPERS-M-9S'
$ws.Range("F371").Value = 6
$ws.Range("G371").Value = 514.08
$ws.Range("B372").Value = 64985
$ws.Range("C372").Value = 'HIM-TOTAL CARE BABY PANTS DRAPERS-XL-9S'
$ws.Range("F372").Value = 13
$ws.Range("G372").Value = 1140.1
$ws.Range("B373").Value = 66196
$ws.Range("C373").Value = 'HIM-Total Care Baby Pants Drapers-Xl-9S'
$ws.Range("F373").Value = 29
$ws.Range("G373").Value = 2543.3
$ws.Range("B483").Value = 47097
$ws.Range("D483").Value = 112.28
$ws.Range("E483").Value = 134.16
$ws.Range("F483").Value = 15
$ws.Range("G483").Value = 1684.2
$ws.Range("B484").Value = 58047
$ws.Range("D484").Value = 105.54
$ws.Range("E484").Value = 126.1
$ws.Range("F484").Value = 34
$ws.Range("G484").Value = 3588.36
$ws.Range("F488").Value = 8
$ws.Range("G488").Value = 2012
$ws.Range("B492").Value = 80588.89999999999
$ws.Range("B551").Value = 53602
$ws.Range("E551").Value = 15.69
$ws.Range("F551").Value = -231
$ws.Range("G551").Value = -3037.65
$ws.Range("B552").Value = 65068
$ws.Range("E552").Value = 13.97
$ws.Range("F552").Value = 63
$ws.Range("G552").Value = 828.45
$ws.Range("B553").Value = 65066
$ws.Range("E553").Value = 13.61
$ws.Range("F553").Value = 90
$ws.Range("G553").Value = 1152.9
$ws.Range("B554").Value = 53263
$ws.Range("E554").Value = 15.29
$ws.Range("F554").Value = -309
$ws.Range("G554").Value = -3958.29
$ws.Range("B559").Value = 45706
$ws.Range("E559").Value = 23.58
$ws.Range("F559").Value = -202
$ws.Range("G559").Value = -3985.46
$ws.Range("B560").Value = 64922
$ws.Range("E560").Value = 20.98
$ws.Range("F560").Value = 67
$ws.Range("G560").Value = 1321.91
$ws.Range("B567").Value = 64925
$ws.Range("E567").Value = 13.97
$ws.Range("F567").Value = 111
$ws.Range("G567").Value = 1459.65
$ws.Range("B568").Value = 45709
$ws.Range("E568").Value = 15.69
$ws.Range("F568").Value = -300
$ws.Range("G568").Value = -3945
$ws.Range("B569").Value = 45702
$ws.Range("E569").Value = 31.43
$ws.Range("F569").Value = -215
$ws.Range("G569").Value = -5654.5
$ws.Range("B570").Value = 64919
$ws.Range("E570").Value = 27.97
$ws.Range("F570").Value = 61
$ws.Range("G570").Value = 1604.3
$ws.Range("F607").Value = 48
$ws.Range("G607").Value = 6780.96
$ws.Range("B615").Value = 149662.55
$ws.Range("F635").Value = 18
$ws.Range("G635").Value = 2795.76
$ws.Range("B636").Value = 72044.98
$ws.Range("B662").Value = 60025
$ws.Range("E662").Value = 37.22
$ws.Range("F662").Value = -98
$ws.Range("G662").Value = -3217.34
$ws.Range("B663").Value = 64833
$ws.Range("E663").Value = 34.9
$ws.Range("F663").Value = 90
$ws.Range("G663").Value = 2954.7
$ws.Range("B672").Value = 64830
$ws.Range("E672").Value = 34.9
$ws.Range("F672").Value = 91
$ws.Range("G672").Value = 2987.53
$ws.Range("B673").Value = 60022
$ws.Range("E673").Value = 37.22
$ws.Range("F673").Value = -113
$ws.Range("G673").Value = -3709.79
$ws.Range("F682").Value = 237
$ws.Range("G682").Value = 11724.39
$ws.Range("F689").Value = 151
$ws.Range("G689").Value = 4137.4
$ws.Range("B695").Value = 187319.08
$ws.Range("F703").Value = 69
$ws.Range("G703").Value = 2979.42
$ws.Range("F705").Value = 65
$ws.Range("G705").Value = 2806.7
$ws.Range("B708").Value = 41444.47
$ws.Range("F753").Value = 12
$ws.Range("G753").Value = 1720.32
$ws.Range("B755").Value = 78247.11
$ws.Range("F805").Value = 10
$ws.Range("G805").Value = 1638.9
$ws.Range("F807").Value = 154
$ws.Range("G807").Value = 16756.74
$ws.Range("F827").Value = 577
$ws.Range("G827").Value = 59379.07
$ws.Range("B838").Value = 331526.26
$ws.Range("F842").Value = 4
$ws.Range("G842").Value = 26950.64
$ws.Range("F844").Value = 4
$ws.Range("G844").Value = 56624.64
$ws.Range("B845").Value = 305885.98
$ws.Range("F862").Value = 2
$ws.Range("G862").Value = 1899.88
$ws.Range("B872").Value = 121925.46
$ws.Range("F890").Value = 239
$ws.Range("G890").Value = 7224.97
$ws.Range("F891").Value = 1816
$ws.Range("G891").Value = 296207.76
$ws.Range("B897").Value = 342647.6
$ws.Range("F931").Value = 140
$ws.Range("G931").Value = 5163.2
$ws.Range("B936").Value = 118257.67
$ws.Range("B942").Value = 5108536.29
$ws.Range("B943").Value = 5108536.29
